# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="47.493.28"; ForceText=$false},
    @{Cell="E2"; Value="  +5.51%  "; ForceText=$false},
    @{Cell="D3"; Value="2.496.26"; ForceText=$false},
    @{Cell="E3"; Value="  +2.89%  "; ForceText=$false},
    @{Cell="E4"; Value="  +0.03%  "; ForceText=$false},
    @{Cell="D5"; Value="323.39"; ForceText=$true},
    @{Cell="E5"; Value="  +2.24%  "; ForceText=$false},
    @{Cell="D6"; Value="105.27"; ForceText=$true},
    @{Cell="E6"; Value="  +2.43%  "; ForceText=$false},
    @{Cell="E7"; Value="  +1.54%  "; ForceText=$false},
    @{Cell="E8"; Value="  +0.02%  "; ForceText=$false},
    @{Cell="E9"; Value="  +2.67%  "; ForceText=$false},
    @{Cell="D10"; Value="37.58"; ForceText=$true},
    @{Cell="E10"; Value="  +6.06%  "; ForceText=$false},
    @{Cell="D11"; Value="0.0813"; ForceText=$true},
    @{Cell="E11"; Value="  +1.37%  "; ForceText=$false},
    @{Cell="E12"; Value="  +0.36%  "; ForceText=$false},
    @{Cell="D13"; Value="18.34"; ForceText=$true},
    @{Cell="E13"; Value="  +0.73%  "; ForceText=$false},
    @{Cell="D14"; Value="7.18"; ForceText=$true},
    @{Cell="E14"; Value="  +2.80%  "; ForceText=$false},
    @{Cell="D15"; Value="2.884.15"; ForceText=$false},
    @{Cell="E15"; Value="  +2.93%  "; ForceText=$false},
    @{Cell="D16"; Value="2.505.25"; ForceText=$false},
    @{Cell="E16"; Value="  +3.64%  "; ForceText=$false},
    @{Cell="E17"; Value="  +0.79%  "; ForceText=$false},
    @{Cell="D18"; Value="47.376.33"; ForceText=$false},
    @{Cell="E18"; Value="  +5.48%  "; ForceText=$false},
    @{Cell="D19"; Value="12.71"; ForceText=$true},
    @{Cell="E19"; Value="  +3.67%  "; ForceText=$false},
    @{Cell="E20"; Value="  +2.91%  "; ForceText=$false},
    @{Cell="D22"; Value="70.78"; ForceText=$true},
    @{Cell="E22"; Value="  +2.93%  "; ForceText=$false},
    @{Cell="D23"; Value="250.61"; ForceText=$true},
    @{Cell="E23"; Value="  +2.87%  "; ForceText=$false},
    @{Cell="E24"; Value="  +5.52%  "; ForceText=$false},
    @{Cell="D25"; Value="2.56"; ForceText=$true},
    @{Cell="E25"; Value="  +2.94%  "; ForceText=$false},
    @{Cell="D26"; Value="26.20"; ForceText=$true},
    @{Cell="E26"; Value="  +3.88%  "; ForceText=$false},
    @{Cell="E27"; Value="  +0.02%  "; ForceText=$false},
    @{Cell="D28"; Value="10.08"; ForceText=$true},
    @{Cell="E28"; Value="  +5.54%  "; ForceText=$false},
    @{Cell="E29"; Value="  +0.86%  "; ForceText=$false},
    @{Cell="D30"; Value="35.20"; ForceText=$true},
    @{Cell="E30"; Value="  +7.24%  "; ForceText=$false},
    @{Cell="E31"; Value="  +8.04%  "; ForceText=$false},
    @{Cell="D32"; Value="49.49"; ForceText=$true},
    @{Cell="E32"; Value="  +0.80%  "; ForceText=$false},
    @{Cell="D33"; Value="20.09"; ForceText=$true},
    @{Cell="E33"; Value="  +0.93%  "; ForceText=$false},
    @{Cell="E34"; Value="  +2.88%  "; ForceText=$false},
    @{Cell="D35"; Value="0.0782"; ForceText=$true},
    @{Cell="E35"; Value="  +2.69%  "; ForceText=$false},
    @{Cell="E36"; Value="  +0.09%  "; ForceText=$false},
    @{Cell="D37"; Value="4.66"; ForceText=$true},
    @{Cell="E37"; Value="  +5.24%  "; ForceText=$false},
    @{Cell="E38"; Value="  +3.75%  "; ForceText=$false},
    @{Cell="E39"; Value="  +4.10%  "; ForceText=$false},
    @{Cell="E40"; Value="  +1.98%  "; ForceText=$false},
    @{Cell="B41"; Value="Monero"; ForceText=$false},
    @{Cell="C41"; Value="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; ForceText=$false},
    @{Cell="D41"; Value="121.05"; ForceText=$true},
    @{Cell="E41"; Value="  -1.36%  "; ForceText=$false},
    @{Cell="B42"; Value="WEMIXToken"; ForceText=$false},
    @{Cell="C42"; Value="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; ForceText=$false},
    @{Cell="D42"; Value="2.23"; ForceText=$true},
    @{Cell="E42"; Value="  +0.79%  "; ForceText=$false},
    @{Cell="D43"; Value="21.15"; ForceText=$true},
    @{Cell="E43"; Value="  +2.24%  "; ForceText=$false},
    @{Cell="D44"; Value="0.0296"; ForceText=$true},
    @{Cell="E44"; Value="  +2.32%  "; ForceText=$false},
    @{Cell="D45"; Value="1.963.91"; ForceText=$false},
    @{Cell="E45"; Value="  +1.53%  "; ForceText=$false},
    @{Cell="E46"; Value="  +1.51%  "; ForceText=$false},
    @{Cell="E47"; Value="  -0.25%  "; ForceText=$false},
    @{Cell="B48"; Value="Stacks"; ForceText=$false},
    @{Cell="C48"; Value="https://coinranking.com/coin/mMPrMcB7+stacks-stx"; ForceText=$false},
    @{Cell="D48"; Value="1.84"; ForceText=$true},
    @{Cell="E48"; Value="  +2.99%  "; ForceText=$false},
    @{Cell="B49"; Value="FraxShare"; ForceText=$false},
    @{Cell="C49"; Value="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; ForceText=$false},
    @{Cell="D49"; Value="9.23"; ForceText=$true},
    @{Cell="E49"; Value="  +0.24%  "; ForceText=$false},
    @{Cell="D50"; Value="5.35"; ForceText=$true},
    @{Cell="E50"; Value="  +14.21%  "; ForceText=$false},
    @{Cell="D51"; Value="78.74"; ForceText=$true},
    @{Cell="E51"; Value="  +3.18%  "; ForceText=$false}
)

foreach ($chg in $changes) {
    $r = $ws.Range($chg.Cell)
    if ($chg.ForceText) {
        $r.NumberFormat = "@"
        $r.Value = $chg.Value
        $r.NumberFormat = "General"
        $r.Style = "Normal"
    } else {
        $r.Value = $chg.Value
    }
}
